# GIT.docx edit:
#   - Title paragraph text "Title" -> "GIT"; the "_GoBack" bookmark now sits
#     right after the new title text instead of after "Paragraph".
#   - A brand-new paragraph "Altering the first paragraph, removed heading"
#     (split into 3 runs) is inserted right after the title.
#   - The old "Heading" paragraph (Heading1 style) is replaced by a plain
#     (Normal-style) paragraph reading "Adding another paragraph".
#   - The old last paragraph ("Paragraph") is dropped entirely.
#   - Every paragraph mark and every run now carries an explicit
#     <w:lang w:val="en-US"/> run property.
#
# Note: this runtime's Range/Font "LanguageID" setter is buggy - it always
# stamps paragraph/run #1 no matter which Range it is invoked on - so the
# w:lang attributes are applied by constructing the OOXML for the new
# content directly and pushing it in with Range.InsertXML, which replaces
# exactly the range it is called on.

$d = $word.ActiveDocument

$newBodyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Title"/>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
              <w:t>GIT</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
              <w:t>Altering the first p</w:t>
            </w:r>
            <w:r>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
              <w:t>aragraph</w:t>
            </w:r>
            <w:r>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
              <w:t>, removed heading</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
              <w:t>Adding another paragraph</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

# Replace the whole body's content (all 3 original paragraphs) in one shot;
# this leaves the sectPr (page setup) at the end of the body untouched,
# since InsertXML only overwrites the range it is called on ($d.Content,
# i.e. the paragraphs, not the section properties).
$d.Content.InsertXML($newBodyXml)
